$d = $word.ActiveDocument
$r = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $r)
Write-Host "done"
